$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) "docker-" / "compose.override" / ".yaml" -> single run
#    "docker-compose.override.yaml" (drop the gramStart/gramEnd proofErr
#    pair that used to straddle the middle run, keep the spellStart/
#    spellEnd pair that wraps the whole token).
# -----------------------------------------------------------------------
$d.Content.Find.Execute("docker-compose.override.yaml", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "docker-compose.override.yaml", 2) | Out-Null

# -----------------------------------------------------------------------
# Journal table for "Mittwoch, 19.06.2024" is the 3rd table in the doc.
# Row 2 = "8:30-10:30", Row 3 = "10:50-12:15".
# Columns: 1 Uhrzeit, 2 Soll, 3 Ist, 4 Probleme.
# -----------------------------------------------------------------------
$t = $d.Tables.Item(3)

# 2) Row "8:30-10:30", "Ist" column (was empty) gets three paragraphs.
$cr = [char]13
$t.Cell(2, 3).Range.Text = "Erstellern der Stored Procedures" + $cr + `
    "Erstellen der API Crud Funktionen " + $cr + `
    "Aktualisieren der Dokumentation"

# 3) Row "8:30-10:30", "Probleme" column (was empty) gets "keine".
$t.Cell(2, 4).Range.Text = "keine"

# 4) Row "10:50-12:15", "Soll" column (was empty) gets two paragraphs.
$t.Cell(3, 2).Range.Text = "Aktualisieren von Zeitplan" + $cr + "Erstellen der Testdaten"
